$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1104.6111
$ws.Range("I53").Value = 1727.3636
$ws.Range("J53").Value = 126
$ws.Range("K53").Value = 1727.3636
$ws.Range("L53").Value = 126
$ws.Range("M53").Value = -1090.3636
$ws.Range("N53").Value = -1400
$ws.Range("H76").Value = 3010
$ws.Range("I76").Value = 3011.111
$ws.Range("K76").Value = 3011.111
$ws.Range("M76").Value = -2696.111
$ws.Range("H79").Value = 3010
$ws.Range("I79").Value = 3011.111
$ws.Range("K79").Value = 3011.111
$ws.Range("M79").Value = -1919.111
$ws.Range("H80").Value = 2384.9375
$ws.Range("J80").Value = 2367.818
$ws.Range("L80").Value = 7103.454000000001
$ws.Range("N80").Value = -9099.454000000002
$ws.Range("H82").Value = 840.1429000000001
$ws.Range("I82").Value = 840.1429000000001
$ws.Range("K82").Value = 2520.4287
$ws.Range("M82").Value = -2114.4287
$ws.Range("H83").Value = 2384.9375
$ws.Range("J83").Value = 2367.818
$ws.Range("L83").Value = 21310.362
$ws.Range("N83").Value = -31294.362
$ws.Range("H85").Value = 840.1429000000001
$ws.Range("I85").Value = 840.1429000000001
$ws.Range("K85").Value = 2520.4287
$ws.Range("M85").Value = -1116.4287
$ws.Range("H88").Value = 6778.7334
$ws.Range("I88").Value = 4980.6
$ws.Range("J88").Value = 7677.8
$ws.Range("K88").Value = 4980.6
$ws.Range("L88").Value = 7677.8
$ws.Range("M88").Value = -4574.6
$ws.Range("N88").Value = -8489.799999999999
$ws.Range("H91").Value = 6778.7334
$ws.Range("I91").Value = 4980.6
$ws.Range("J91").Value = 7677.8
$ws.Range("K91").Value = 4980.6
$ws.Range("L91").Value = 7677.8
$ws.Range("M91").Value = -3576.6
$ws.Range("N91").Value = -10485.8
$ws.Range("H137").Value = 2042.5
$ws.Range("I137").Value = 1490.909
$ws.Range("K137").Value = 4472.727000000001
$ws.Range("M137").Value = -1922.727000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2023.9286
$ws.Range("J88").Value = 1920.375
$ws.Range("L88").Value = 1920.375
$ws.Range("N88").Value = -2732.375
$ws.Range("H91").Value = 2023.9286
$ws.Range("J91").Value = 1920.375
$ws.Range("L91").Value = 1920.375
$ws.Range("N91").Value = -4728.375

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1639.7333
$ws.Range("I86").Value = 1440.5454
$ws.Range("J86").Value = 2187.5
$ws.Range("K86").Value = 1440.5454
$ws.Range("L86").Value = 2187.5
$ws.Range("M86").Value = -317.5454
$ws.Range("N86").Value = -4433.5
$ws.Range("H89").Value = 1639.7333
$ws.Range("I89").Value = 1440.5454
$ws.Range("J89").Value = 2187.5
$ws.Range("K89").Value = 7202.727
$ws.Range("L89").Value = 10937.5
$ws.Range("M89").Value = -1586.727
$ws.Range("N89").Value = -22169.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1445.6066
$ws.Range("I31").Value = 896.625
$ws.Range("K31").Value = 896.625
$ws.Range("M31").Value = -601.625
$ws.Range("H34").Value = 1445.6066
$ws.Range("I34").Value = 896.625
$ws.Range("K34").Value = 896.625
$ws.Range("M34").Value = -694.625
$ws.Range("H62").Value = 2355
$ws.Range("J62").Value = 2450
$ws.Range("L62").Value = 2450
$ws.Range("N62").Value = -3698
$ws.Range("H65").Value = 2355
$ws.Range("J65").Value = 2450
$ws.Range("L65").Value = 12250
$ws.Range("N65").Value = -18490
$ws.Range("H132").Value = 927890.0600000001
$ws.Range("I132").Value = 1939.625
$ws.Range("J132").Value = 4631692
$ws.Range("K132").Value = 5818.875
$ws.Range("L132").Value = 13895076
$ws.Range("M132").Value = -3288.875
$ws.Range("N132").Value = -13900136

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 914.25
$ws.Range("I131").Value = 530
$ws.Range("J131").Value = 918.1313
$ws.Range("K131").Value = 1590
$ws.Range("L131").Value = 2754.3939
$ws.Range("M131").Value = 3450
$ws.Range("N131").Value = -12834.3939

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2700.7144
$ws.Range("I80").Value = 2321
$ws.Range("J80").Value = 3650
$ws.Range("K80").Value = 2321
$ws.Range("L80").Value = 3650
$ws.Range("M80").Value = -1323
$ws.Range("N80").Value = -5646
$ws.Range("H83").Value = 2700.7144
$ws.Range("I83").Value = 2321
$ws.Range("J83").Value = 3650
$ws.Range("K83").Value = 11605
$ws.Range("L83").Value = 18250
$ws.Range("M83").Value = -6613
$ws.Range("N83").Value = -28234
$ws.Range("H132").Value = 2634654.8
$ws.Range("I132").Value = 3374.6956
$ws.Range("J132").Value = 6669284.5
$ws.Range("K132").Value = 10124.0868
$ws.Range("L132").Value = 20007853.5
$ws.Range("M132").Value = -7594.086800000001
$ws.Range("N132").Value = -20012913.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 515.9167
$ws.Range("J55").Value = 743.4286
$ws.Range("L55").Value = 743.4286
$ws.Range("N55").Value = -1089.4286
$ws.Range("H68").Value = 13090
$ws.Range("I68").Value = 35000
$ws.Range("J68").Value = 3700
$ws.Range("K68").Value = 35000
$ws.Range("L68").Value = 3700
$ws.Range("M68").Value = -34251
$ws.Range("N68").Value = -5198
$ws.Range("H71").Value = 13090
$ws.Range("I71").Value = 35000
$ws.Range("J71").Value = 3700
$ws.Range("K71").Value = 175000
$ws.Range("L71").Value = 18500
$ws.Range("M71").Value = -171256
$ws.Range("N71").Value = -25988
$ws.Range("H82").Value = 1446.8572
$ws.Range("I82").Value = 1810.5
$ws.Range("K82").Value = 1810.5
$ws.Range("M82").Value = -1449.5
$ws.Range("H85").Value = 1446.8572
$ws.Range("I85").Value = 1810.5
$ws.Range("K85").Value = 1810.5
$ws.Range("M85").Value = -562.5
$ws.Range("H132").Value = 3244.0889
$ws.Range("I132").Value = 3260.162
$ws.Range("J132").Value = 3169.75
$ws.Range("K132").Value = 9780.485999999999
$ws.Range("L132").Value = 9509.25
$ws.Range("M132").Value = -7250.485999999999
$ws.Range("N132").Value = -14569.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 3500
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 3500
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -2876
$ws.Range("N62").Value = -3248
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 3500
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 17500
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -14380
$ws.Range("N65").Value = -16240
$ws.Range("H81").Value = 1206.1177
$ws.Range("J81").Value = 1220
$ws.Range("L81").Value = 2440
$ws.Range("N81").Value = -4562
$ws.Range("H84").Value = 1206.1177
$ws.Range("J84").Value = 1220
$ws.Range("N84").Value = -22808
